$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-05-20 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-21 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("431×8=3448", $true, $false, $false, $false, $false, $true, 1, $false, "387×4=1548", 2) | Out-Null
$d.Content.Find.Execute("153×3=459", $true, $false, $false, $false, $false, $true, 1, $false, "905×7=6335", 2) | Out-Null
$d.Content.Find.Execute("315×6=1890", $true, $false, $false, $false, $false, $true, 1, $false, "861×4=3444", 2) | Out-Null
$d.Content.Find.Execute("716×2=1432", $true, $false, $false, $false, $false, $true, 1, $false, "895×5=4475", 2) | Out-Null
$d.Content.Find.Execute("443×6=2658", $true, $false, $false, $false, $false, $true, 1, $false, "504×3=1512", 2) | Out-Null
$d.Content.Find.Execute("379×7=2653", $true, $false, $false, $false, $false, $true, 1, $false, "396×8=3168", 2) | Out-Null
$d.Content.Find.Execute("932×9=8388", $true, $false, $false, $false, $false, $true, 1, $false, "230×8=1840", 2) | Out-Null
$d.Content.Find.Execute("982×6=5892", $true, $false, $false, $false, $false, $true, 1, $false, "622×9=5598", 2) | Out-Null
$d.Content.Find.Execute("685×8=5480", $true, $false, $false, $false, $false, $true, 1, $false, "236×7=1652", 2) | Out-Null
$d.Content.Find.Execute("465×7=3255", $true, $false, $false, $false, $false, $true, 1, $false, "259×5=1295", 2) | Out-Null
$d.Content.Find.Execute("616×4=2464", $true, $false, $false, $false, $false, $true, 1, $false, "662×3=1986", 2) | Out-Null
$d.Content.Find.Execute("754×3=2262", $true, $false, $false, $false, $false, $true, 1, $false, "219×5=1095", 2) | Out-Null
$d.Content.Find.Execute("862×2=1724", $true, $false, $false, $false, $false, $true, 1, $false, "218×4=872", 2) | Out-Null
$d.Content.Find.Execute("333×8=2664", $true, $false, $false, $false, $false, $true, 1, $false, "638×4=2552", 2) | Out-Null
$d.Content.Find.Execute("944×3=2832", $true, $false, $false, $false, $false, $true, 1, $false, "313×3=939", 2) | Out-Null
$d.Content.Find.Execute("222×9=1998", $true, $false, $false, $false, $false, $true, 1, $false, "491×7=3437", 2) | Out-Null
$d.Content.Find.Execute("990×2=1980", $true, $false, $false, $false, $false, $true, 1, $false, "780×7=5460", 2) | Out-Null
$d.Content.Find.Execute("807×4=3228", $true, $false, $false, $false, $false, $true, 1, $false, "244×2=488", 2) | Out-Null
$d.Content.Find.Execute("190×5=950", $true, $false, $false, $false, $false, $true, 1, $false, "274×7=1918", 2) | Out-Null
$d.Content.Find.Execute("757×5=3785", $true, $false, $false, $false, $false, $true, 1, $false, "696×7=4872", 2) | Out-Null
$d.Content.Find.Execute("166×2=332", $true, $false, $false, $false, $false, $true, 1, $false, "884×7=6188", 2) | Out-Null
$d.Content.Find.Execute("722×6=4332", $true, $false, $false, $false, $false, $true, 1, $false, "595×7=4165", 2) | Out-Null
$d.Content.Find.Execute("642×6=3852", $true, $false, $false, $false, $false, $true, 1, $false, "765×8=6120", 2) | Out-Null
$d.Content.Find.Execute("932×2=1864", $true, $false, $false, $false, $false, $true, 1, $false, "709×3=2127", 2) | Out-Null
$d.Content.Find.Execute("390×7=2730", $true, $false, $false, $false, $false, $true, 1, $false, "496×9=4464", 2) | Out-Null
